$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 12 have their "Förändrad" date (column C) bumped by one day,
# from serial 45243 to 45244.
for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
